$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.646.98'
$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").Value = '3.425.27'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.40'
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("E6").Value = '  -2.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  +7.12%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '3.429.23'
$ws.Range("E9").Value = '  -0.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.17'
$ws.Range("E10").Value = '  -2.62%  '

$ws.Range("E11").Value = '  -1.51%  '

$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").Value = '4.017.39'
$ws.Range("E13").Value = '  -1.05%  '

$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("E15").Value = '  -3.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.91'
$ws.Range("E16").Value = '  -0.91%  '

$ws.Range("D17").Value = '64.628.74'
$ws.Range("E17").Value = '  -0.46%  '

$ws.Range("D18").Value = '3.450.90'
$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("E20").Value = '  -2.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '378.80'
$ws.Range("E21").Value = '  -2.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.05'
$ws.Range("E22").Value = '  -1.96%  '

$ws.Range("E23").Value = '  +1.17%  '

$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.12%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.53'
$ws.Range("E25").Value = '  -0.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  -3.71%  '

$ws.Range("E27").Value = '  +6.02%  '

$ws.Range("E28").Value = '  -1.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.16%  '

$ws.Range("E30").Value = '  +4.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.22'
$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("E32").Value = '  -0.59%  '

$ws.Range("E33").Value = '  -2.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.23'
$ws.Range("E34").Value = '  +2.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  +7.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.61'
$ws.Range("E36").Value = '  -2.15%  '

$ws.Range("E37").Value = '  -1.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.99'
$ws.Range("E38").Value = '  +6.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0765'
$ws.Range("E39").Value = '  -0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.00'
$ws.Range("E40").Value = '  -0.96%  '

$ws.Range("D41").Value = '2.877.69'
$ws.Range("E41").Value = '  -4.21%  '

$ws.Range("E42").Value = '  +1.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.63'
$ws.Range("E43").Value = '  +9.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0319'
$ws.Range("E44").Value = '  +0.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.92'
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.774'
$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.01'
$ws.Range("E47").Value = '  +5.20%  '

$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.20'
$ws.Range("E50").Value = '  +1.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.59'
$ws.Range("E51").Value = '  -0.06%  '
